# Leave Card update - "Add Leave Card 10/32023 3:18 PM"
# Applies the leave-credit entries for row 60 (correction) and rows 63-71
# (VL/SL/SP earned entries) on the "2018 LEAVE CREDITS" sheet of
# Juel Coper's leave card workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# --- Row 60 (period ending 1/23/2023): was SL(1-0-0) w/ Absence Undertime
#     W/O Pay flag, corrected to VL(1-0-0) w/ Absence Undertime W/ Pay flag,
#     and the remarks date moves from 1/23/2023 to 1/30/2023.
$ws.Range("B60").Value = "VL(1-0-0)"
$ws.Range("D60").Value = 1
$ws.Range("H60").ClearContents()
$ws.Range("K60").Value = 44956

# --- Row 63 (period ending 4/30/2023): flag VL(1-0-0) used, with the
#     Absence Undertime W/ Pay checkbox.
$ws.Range("B63").Value = "VL(1-0-0)"
$ws.Range("D63").Value = 1

# --- Row 64 (period ending 5/31/2023): same as row 63.
$ws.Range("B64").Value = "VL(1-0-0)"
$ws.Range("D64").Value = 1

# --- Row 65 (period ending 6/30/2023): SP(1-0-0) earned 1.25, remarks
#     dated 6/15/2023. The remarks cell was blank (no number format), so
#     entering a date makes Excel pick up the existing "date" style used
#     by the other remarks-date cells (copy format from K60).
$ws.Range("B65").Value = "SP(1-0-0)"
$ws.Range("C65").Value = 1.25
$ws.Range("K60").Copy() | Out-Null
$ws.Range("K65").PasteSpecial(-4122) | Out-Null
$ws.Range("K65").Value = 45092

# --- Row 66 (period ending 7/31/2023): SL(1-0-0) earned 1.25, Absence
#     Undertime W/ Pay flagged, remarks dated 7/26/2023.
$ws.Range("B66").Value = "SL(1-0-0)"
$ws.Range("C66").Value = 1.25
$ws.Range("H66").Value = 1
$ws.Range("K60").Copy() | Out-Null
$ws.Range("K66").PasteSpecial(-4122) | Out-Null
$ws.Range("K66").Value = 45133

# --- Row 67 (period ending 8/31/2023): VL(1-0-0) earned 1.25, Absence
#     Undertime W/ Pay flagged, remarks dated 9/1/2023.
$ws.Range("B67").Value = "VL(1-0-0)"
$ws.Range("C67").Value = 1.25
$ws.Range("D67").Value = 1
$ws.Range("K60").Copy() | Out-Null
$ws.Range("K67").PasteSpecial(-4122) | Out-Null
$ws.Range("K67").Value = 45170

$excel.CutCopyMode = $false

# --- Rows 68-71 (periods ending 9/30, 10/31, 11/30, 12/31/2023): plain
#     1.25 leave credit earned each month, no particulars/remarks.
$ws.Range("C68").Value = 1.25
$ws.Range("C69").Value = 1.25
$ws.Range("C70").Value = 1.25
$ws.Range("C71").Value = 1.25
